# Apply the "leak" columns removal + RES relabeling described by the commit
# "changes w.r.t res and leaks":
#   - Columns Z1:AC1 (currently the first 4 "...Leak" headers) become the
#     RES2 / RES4 / RES16 / RES20 headers that used to live at the far end
#     of the table (columns BA:BD).
#   - All the remaining "...Leak" columns (AD:BD) are removed entirely,
#     shrinking the used range from A1:BD2 down to A1:AC2.
#   - The leak-rate sample values that used to sit under Z:AC are replaced
#     with 0 (matching the 0s that already filled the now-deleted BA:BD
#     columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-point the first four "Leak" header columns to the RES* labels.
$ws.Range("Z1").Value  = "RES2"
$ws.Range("AA1").Value = "RES4"
$ws.Range("AB1").Value = "RES16"
$ws.Range("AC1").Value = "RES20"

# Zero out the corresponding data row for those columns.
$ws.Range("Z2:AC2").Value = 0

# Drop all the remaining leak columns (and the old trailing RES* columns,
# now duplicated above) so the sheet shrinks to A1:AC2.
$ws.Range("AD1:BD2").EntireColumn.Delete()
